# ============================================================================
# edit.ps1 - Apply "Update test results data" commit changes to
# DU FBS Mock 7.xlsx (Sheet3 / "Sheet3") using Excel COM-interop semantics.
#
# Summary of the edit:
#   1. Workbook-level: Sheet3 becomes the active tab (Sheet1 loses it).
#   2. Sheet3's sheetView: bottomRight pane selection becomes A22:XFD22.
#   3. Sheet3 row 4 (columns C:BZ) and row 22 (columns C:BZ) are populated
#      with quiz-answer-key values ("A (C)", "B (W)", ... ) reusing the
#      workbook's existing shared strings, matching the style already used
#      by sibling rows (4/5/23 etc.) in the same columns.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# Step 1: Data for rows 4 and 22 (columns C through BZ).
#   S 16/17 cells carry one of the 8 recurring grade strings; S 18 cells are
#   intentionally blank (no value), matching neighbouring data rows already
#   in the sheet (e.g. row 5, row 23).
# ---------------------------------------------------------------------------
$row4Data = @(
    @{ Col = "C"; S = 16; Val = "B (C)" }
    @{ Col = "D"; S = 16; Val = "A (C)" }
    @{ Col = "E"; S = 17; Val = "A (W)" }
    @{ Col = "F"; S = 16; Val = "D (C)" }
    @{ Col = "G"; S = 17; Val = "A (W)" }
    @{ Col = "H"; S = 16; Val = "B (C)" }
    @{ Col = "I"; S = 16; Val = "A (C)" }
    @{ Col = "J"; S = 16; Val = "B (C)" }
    @{ Col = "K"; S = 16; Val = "B (C)" }
    @{ Col = "L"; S = 16; Val = "A (C)" }
    @{ Col = "M"; S = 16; Val = "B (C)" }
    @{ Col = "N"; S = 16; Val = "C (C)" }
    @{ Col = "O"; S = 16; Val = "B (C)" }
    @{ Col = "P"; S = 16; Val = "D (C)" }
    @{ Col = "Q"; S = 17; Val = "A (W)" }
    @{ Col = "R"; S = 17; Val = "B (W)" }
    @{ Col = "S"; S = 16; Val = "C (C)" }
    @{ Col = "T"; S = 17; Val = "B (W)" }
    @{ Col = "U"; S = 16; Val = "B (C)" }
    @{ Col = "V"; S = 16; Val = "A (C)" }
    @{ Col = "W"; S = 17; Val = "D (W)" }
    @{ Col = "X"; S = 17; Val = "D (W)" }
    @{ Col = "Y"; S = 17; Val = "C (W)" }
    @{ Col = "Z"; S = 16; Val = "A (C)" }
    @{ Col = "AA"; S = 16; Val = "C (C)" }
    @{ Col = "AB"; S = 17; Val = "C (W)" }
    @{ Col = "AC"; S = 17; Val = "D (W)" }
    @{ Col = "AD"; S = 16; Val = "C (C)" }
    @{ Col = "AE"; S = 17; Val = "B (W)" }
    @{ Col = "AF"; S = 16; Val = "C (C)" }
    @{ Col = "AG"; S = 16; Val = "B (C)" }
    @{ Col = "AH"; S = 16; Val = "C (C)" }
    @{ Col = "AI"; S = 16; Val = "D (C)" }
    @{ Col = "AJ"; S = 16; Val = "C (C)" }
    @{ Col = "AK"; S = 16; Val = "D (C)" }
    @{ Col = "AL"; S = 16; Val = "C (C)" }
    @{ Col = "AM"; S = 16; Val = "A (C)" }
    @{ Col = "AN"; S = 16; Val = "C (C)" }
    @{ Col = "AO"; S = 16; Val = "D (C)" }
    @{ Col = "AP"; S = 17; Val = "C (W)" }
    @{ Col = "AQ"; S = 16; Val = "D (C)" }
    @{ Col = "AR"; S = 17; Val = "B (W)" }
    @{ Col = "AS"; S = 17; Val = "C (W)" }
    @{ Col = "AT"; S = 17; Val = "A (W)" }
    @{ Col = "AU"; S = 18; Val = $null }
    @{ Col = "AV"; S = 18; Val = $null }
    @{ Col = "AW"; S = 18; Val = $null }
    @{ Col = "AX"; S = 18; Val = $null }
    @{ Col = "AY"; S = 18; Val = $null }
    @{ Col = "AZ"; S = 18; Val = $null }
    @{ Col = "BA"; S = 18; Val = $null }
    @{ Col = "BB"; S = 18; Val = $null }
    @{ Col = "BC"; S = 18; Val = $null }
    @{ Col = "BD"; S = 18; Val = $null }
    @{ Col = "BE"; S = 18; Val = $null }
    @{ Col = "BF"; S = 18; Val = $null }
    @{ Col = "BG"; S = 18; Val = $null }
    @{ Col = "BH"; S = 18; Val = $null }
    @{ Col = "BI"; S = 18; Val = $null }
    @{ Col = "BJ"; S = 18; Val = $null }
    @{ Col = "BK"; S = 16; Val = "C (C)" }
    @{ Col = "BL"; S = 16; Val = "A (C)" }
    @{ Col = "BM"; S = 16; Val = "C (C)" }
    @{ Col = "BN"; S = 17; Val = "D (W)" }
    @{ Col = "BO"; S = 16; Val = "C (C)" }
    @{ Col = "BP"; S = 16; Val = "B (C)" }
    @{ Col = "BQ"; S = 16; Val = "B (C)" }
    @{ Col = "BR"; S = 16; Val = "C (C)" }
    @{ Col = "BS"; S = 17; Val = "D (W)" }
    @{ Col = "BT"; S = 16; Val = "B (C)" }
    @{ Col = "BU"; S = 16; Val = "C (C)" }
    @{ Col = "BV"; S = 16; Val = "B (C)" }
    @{ Col = "BW"; S = 16; Val = "B (C)" }
    @{ Col = "BX"; S = 18; Val = $null }
    @{ Col = "BY"; S = 16; Val = "A (C)" }
    @{ Col = "BZ"; S = 16; Val = "C (C)" }
)

$row22Data = @(
    @{ Col = "C"; S = 16; Val = "B (C)" }
    @{ Col = "D"; S = 16; Val = "A (C)" }
    @{ Col = "E"; S = 17; Val = "A (W)" }
    @{ Col = "F"; S = 16; Val = "D (C)" }
    @{ Col = "G"; S = 16; Val = "C (C)" }
    @{ Col = "H"; S = 16; Val = "B (C)" }
    @{ Col = "I"; S = 16; Val = "A (C)" }
    @{ Col = "J"; S = 16; Val = "B (C)" }
    @{ Col = "K"; S = 16; Val = "B (C)" }
    @{ Col = "L"; S = 16; Val = "A (C)" }
    @{ Col = "M"; S = 16; Val = "B (C)" }
    @{ Col = "N"; S = 16; Val = "C (C)" }
    @{ Col = "O"; S = 16; Val = "B (C)" }
    @{ Col = "P"; S = 16; Val = "D (C)" }
    @{ Col = "Q"; S = 16; Val = "B (C)" }
    @{ Col = "R"; S = 17; Val = "B (W)" }
    @{ Col = "S"; S = 16; Val = "C (C)" }
    @{ Col = "T"; S = 17; Val = "D (W)" }
    @{ Col = "U"; S = 16; Val = "B (C)" }
    @{ Col = "V"; S = 16; Val = "A (C)" }
    @{ Col = "W"; S = 16; Val = "B (C)" }
    @{ Col = "X"; S = 16; Val = "B (C)" }
    @{ Col = "Y"; S = 18; Val = $null }
    @{ Col = "Z"; S = 16; Val = "A (C)" }
    @{ Col = "AA"; S = 16; Val = "C (C)" }
    @{ Col = "AB"; S = 16; Val = "B (C)" }
    @{ Col = "AC"; S = 17; Val = "D (W)" }
    @{ Col = "AD"; S = 16; Val = "C (C)" }
    @{ Col = "AE"; S = 16; Val = "D (C)" }
    @{ Col = "AF"; S = 16; Val = "C (C)" }
    @{ Col = "AG"; S = 16; Val = "B (C)" }
    @{ Col = "AH"; S = 16; Val = "C (C)" }
    @{ Col = "AI"; S = 17; Val = "C (W)" }
    @{ Col = "AJ"; S = 16; Val = "C (C)" }
    @{ Col = "AK"; S = 16; Val = "D (C)" }
    @{ Col = "AL"; S = 16; Val = "C (C)" }
    @{ Col = "AM"; S = 17; Val = "D (W)" }
    @{ Col = "AN"; S = 16; Val = "C (C)" }
    @{ Col = "AO"; S = 16; Val = "D (C)" }
    @{ Col = "AP"; S = 16; Val = "B (C)" }
    @{ Col = "AQ"; S = 16; Val = "D (C)" }
    @{ Col = "AR"; S = 16; Val = "D (C)" }
    @{ Col = "AS"; S = 16; Val = "B (C)" }
    @{ Col = "AT"; S = 16; Val = "C (C)" }
    @{ Col = "AU"; S = 18; Val = $null }
    @{ Col = "AV"; S = 18; Val = $null }
    @{ Col = "AW"; S = 18; Val = $null }
    @{ Col = "AX"; S = 18; Val = $null }
    @{ Col = "AY"; S = 18; Val = $null }
    @{ Col = "AZ"; S = 18; Val = $null }
    @{ Col = "BA"; S = 18; Val = $null }
    @{ Col = "BB"; S = 18; Val = $null }
    @{ Col = "BC"; S = 18; Val = $null }
    @{ Col = "BD"; S = 18; Val = $null }
    @{ Col = "BE"; S = 18; Val = $null }
    @{ Col = "BF"; S = 18; Val = $null }
    @{ Col = "BG"; S = 18; Val = $null }
    @{ Col = "BH"; S = 18; Val = $null }
    @{ Col = "BI"; S = 18; Val = $null }
    @{ Col = "BJ"; S = 18; Val = $null }
    @{ Col = "BK"; S = 17; Val = "B (W)" }
    @{ Col = "BL"; S = 17; Val = "C (W)" }
    @{ Col = "BM"; S = 16; Val = "C (C)" }
    @{ Col = "BN"; S = 17; Val = "D (W)" }
    @{ Col = "BO"; S = 16; Val = "C (C)" }
    @{ Col = "BP"; S = 16; Val = "B (C)" }
    @{ Col = "BQ"; S = 16; Val = "B (C)" }
    @{ Col = "BR"; S = 16; Val = "C (C)" }
    @{ Col = "BS"; S = 17; Val = "C (W)" }
    @{ Col = "BT"; S = 18; Val = $null }
    @{ Col = "BU"; S = 17; Val = "A (W)" }
    @{ Col = "BV"; S = 16; Val = "B (C)" }
    @{ Col = "BW"; S = 17; Val = "D (W)" }
    @{ Col = "BX"; S = 16; Val = "C (C)" }
    @{ Col = "BY"; S = 16; Val = "A (C)" }
    @{ Col = "BZ"; S = 16; Val = "C (C)" }
)

# Donor cells already carrying styles 16 / 17 / 18 in the sheet, used purely
# as a formatting source (Copy + PasteSpecial formats) so the new cells end
# up on the exact same cellXfs index as their neighbours.
$donor16 = $ws3.Range("C5")
$donor17 = $ws3.Range("E5")
$donor18 = $ws3.Range("H5")

function Apply-RowData($ws, $rowNum, $rowData) {
    foreach ($entry in $rowData) {
        $target = $ws.Range($entry.Col + $rowNum)

        switch ($entry.S) {
            16 { $donor16.Copy(); $target.PasteSpecial(-4122) }
            17 { $donor17.Copy(); $target.PasteSpecial(-4122) }
            18 { $donor18.Copy(); $target.PasteSpecial(-4122) }
        }

        if ($null -ne $entry.Val) {
            $target.Value = $entry.Val
        }
    }
}

Apply-RowData $ws3 4 $row4Data
Apply-RowData $ws3 22 $row22Data

# ---------------------------------------------------------------------------
# Step 2: Switch the active/selected sheet from Sheet1 to Sheet3, and set
# Sheet3's frozen-pane ("bottomRight") selection to A22:XFD22 (full row 22),
# matching the saved selection captured in the workbook after editing.
# ---------------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("A22:XFD22").Select()

Write-Host "Edit complete: Sheet3 rows 4 & 22 populated; Sheet3 is now the active tab."
